$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("D8").Value = "AC0603JRNPO9BN220"
$ws.Range("E8").Value = "603-AC603JRNPO9BN220"
